$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.697.76"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.644.34"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.530"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.876.62"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.647.40"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "27.666.27"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.34%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "1.444.37"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.570"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.881"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.891"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.83%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").Value = "1.786.52"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("E47").Value = "  +6.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0988"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0504"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
